$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 used to repeat the same loan record as rows 3/5 (StB bank, "МКД"),
# but carried a different long note in L2 (the "Максимум возраст ... 5.55%"
# text, which was unique to this row). The row is being stripped back down
# to a near-blank placeholder: only the bank name stays in column A, while
# H2 and L2 keep their number/wrap-text styling but lose their values, and
# every other cell in the row is removed outright.
$ws.Range("B2:G2").Clear()
$ws.Range("I2:K2").Clear()
$ws.Range("M2").Clear()
$ws.Range("H2").ClearContents()
$ws.Range("L2").ClearContents()

# With the long note gone, row 2 no longer needs its tall custom height -
# let it size back down to the sheet's normal row height.
$ws.Rows("2").AutoFit()
